$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings (like "41.17")
# are preserved as text instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.808.21'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").Value = '1.841.67'
$ws.Range("E3").Value = '  -1.54%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").Value = '308.56'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").Value = '0.9938'
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("D7").Value = '0.5029'
$ws.Range("E7").Value = '  -2.51%  '
$ws.Range("D8").Value = '0.3847'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.08106'
$ws.Range("E9").Value = '  -2.07%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").Value = '41.17'
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '1.096'
$ws.Range("E11").Value = '  -1.37%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '6.103'
$ws.Range("E12").Value = '  -1.72%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.820.87'
$ws.Range("E13").Value = '  -3.22%  '
$ws.Range("B14").Value = 'Solana'
$ws.Range("C14").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D14").Value = '19.95'
$ws.Range("E14").Value = '  -2.96%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '7.108'
$ws.Range("E15").Value = '  -2.73%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '0.9977'
$ws.Range("E16").Value = '  -0.61%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.00001085'
$ws.Range("E17").Value = '  -1.25%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '89.95'
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.06599'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '17.42'
$ws.Range("E20").Value = '  -1.62%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '0.9923'
$ws.Range("E21").Value = '  -1.04%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '5.894'
$ws.Range("E22").Value = '  -2.33%  '
$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D23").Value = '27.846.68'
$ws.Range("E23").Value = '  -0.77%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '10.92'
$ws.Range("E24").Value = '  -1.68%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.202'
$ws.Range("E25").Value = '  -2.25%  '
$ws.Range("D26").Value = '2.046.74'
$ws.Range("E26").Value = '  -1.85%  '
$ws.Range("D27").Value = '157.53'
$ws.Range("E27").Value = '  +0.07%  '
$ws.Range("D28").Value = '20.39'
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("D29").Value = '2.373'
$ws.Range("E29").Value = '  -5.10%  '
$ws.Range("D30").Value = '124.50'
$ws.Range("E30").Value = '  -0.54%  '
$ws.Range("D31").Value = '0.1043'
$ws.Range("E31").Value = '  -2.29%  '
$ws.Range("D32").Value = '1.026'
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("D33").Value = '5.754'
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("D34").Value = '3.566'
$ws.Range("E34").Value = '  -0.81%  '
$ws.Range("D35").Value = '0.02390'
$ws.Range("E35").Value = '  -1.44%  '
$ws.Range("D36").Value = '9.121'
$ws.Range("E36").Value = '  -3.94%  '
$ws.Range("D37").Value = '0.06434'
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("D38").Value = '0.2155'
$ws.Range("E38").Value = '  -2.24%  '
$ws.Range("D39").Value = '0.6329'
$ws.Range("E39").Value = '  -3.80%  '
$ws.Range("D40").Value = '1.217'
$ws.Range("E40").Value = '  +0.60%  '
$ws.Range("D41").Value = '1.158'
$ws.Range("E41").Value = '  -3.69%  '
$ws.Range("D42").Value = '4.887'
$ws.Range("E42").Value = '  -2.60%  '
$ws.Range("D43").Value = '11.00'
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("D44").Value = '0.5965'
$ws.Range("E44").Value = '  -2.96%  '
$ws.Range("D45").Value = '12.92'
$ws.Range("E45").Value = '  -1.68%  '
$ws.Range("D46").Value = '1.261'
$ws.Range("E46").Value = '  -1.84%  '
$ws.Range("D47").Value = '3.622'
$ws.Range("E47").Value = '  -1.17%  '
$ws.Range("D48").Value = '1.964'
$ws.Range("E48").Value = '  -3.19%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '1.187'
$ws.Range("E49").Value = '  -2.71%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = '119.93'
$ws.Range("E50").Value = '  -0.77%  '
$ws.Range("D51").Value = '0.06827'
$ws.Range("E51").Value = '  +0.09%  '
# Restore default styling (no explicit style index) on the edited range.
$dataRange.Style = "Normal"
